$d = $word.ActiveDocument

$replacements = @(
    @("94÷2=", "40÷4="),
    @("37÷5=", "66÷7="),
    @("77÷9=", "54÷3="),
    @("19÷3=", "11÷4="),
    @("31÷6=", "43÷3="),
    @("29÷9=", "72÷7="),
    @("17÷3=", "75÷7="),
    @("73÷3=", "68÷4="),
    @("10÷6=", "34÷8="),
    @("43÷6=", "18÷8="),
    @("34÷4=", "44÷2="),
    @("63÷3=", "35÷3="),
    @("59÷3=", "67÷2="),
    @("60÷5=", "18÷3="),
    @("80÷9=", "99÷2="),
    @("83÷9=", "88÷9="),
    @("49÷5=", "91÷4="),
    @("40÷9=", "62÷7="),
    @("63÷6=", "32÷7="),
    @("90÷8=", "11÷2="),
    @("25÷8=", "59÷2="),
    @("33÷3=", "25÷7="),
    @("74÷7=", "17÷4="),
    @("12÷2=", "11÷6="),
    @("64÷7=", "24÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
